$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the data rows of columns D and E as Text before writing so values
# keep their exact string representation (e.g. trailing zeros, percent
# signs) instead of being auto-converted into numbers/percentages by Excel.
# (Row 1 is the header and is intentionally excluded so its style is kept.)
$dataRange = "D2:E51"
$ws.Range($dataRange).NumberFormat = "@"

$ws.Range("D2").Value = "299.78"
$ws.Range("E2").Value = "-0.35%"
$ws.Range("D3").Value = "31.82"
$ws.Range("E3").Value = "1.50%"
$ws.Range("D4").Value = "5.145"
$ws.Range("E4").Value = "0.60%"
$ws.Range("D5").Value = "0.08187"
$ws.Range("E5").Value = "11.25%"
$ws.Range("D6").Value = "2.550"
$ws.Range("E6").Value = "9.84%"
$ws.Range("D7").Value = "7.852"
$ws.Range("E7").Value = "-1.21%"
$ws.Range("D8").Value = "3.864"
$ws.Range("E8").Value = "2.05%"
$ws.Range("D9").Value = "0.9270"
$ws.Range("E9").Value = "0.93%"
$ws.Range("D10").Value = "0.1758"
$ws.Range("D11").Value = "0.07487"
$ws.Range("E11").Value = "-1.67%"
$ws.Range("D12").Value = "0.08969"
$ws.Range("E12").Value = "10.71%"
$ws.Range("D13").Value = "0.03027"
$ws.Range("E13").Value = "1.41%"
$ws.Range("D14").Value = "0.1003"
$ws.Range("E14").Value = "1.07%"
$ws.Range("D15").Value = "0.001511"
$ws.Range("E15").Value = "1.10%"
$ws.Range("D16").Value = "0.006042"
$ws.Range("E16").Value = "-2.28%"
$ws.Range("D17").Value = "3.603"
$ws.Range("E17").Value = "4.01%"
$ws.Range("E19").Value = "-0.92%"
$ws.Range("D20").Value = "0.1347"
$ws.Range("E20").Value = "2.16%"
$ws.Range("D21").Value = "4.246"
$ws.Range("E21").Value = "-8.80%"
$ws.Range("D22").Value = "0.1678"
$ws.Range("E22").Value = "7.12%"
$ws.Range("D23").Value = "0.04635"
$ws.Range("E23").Value = "0.06%"
$ws.Range("D24").Value = "0.001248"
$ws.Range("E24").Value = "1.83%"
$ws.Range("D25").Value = "0.004551"
$ws.Range("E25").Value = "1.61%"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").Value = "-7.54%"
$ws.Range("D27").Value = "0.0003405"
$ws.Range("E27").Value = "81.95%"
$ws.Range("E39").Value = "3.11%"
$ws.Range("D40").Value = "0.04604"
$ws.Range("E40").Value = "2.10%"
$ws.Range("D41").Value = "0.006864"
$ws.Range("E41").Value = "-5.26%"
$ws.Range("D42").Value = "0.1382"
$ws.Range("E42").Value = "2.85%"
$ws.Range("D43").Value = "0.002142"
$ws.Range("E43").Value = "-4.29%"
$ws.Range("D44").Value = "0.009864"
$ws.Range("E44").Value = "-7.98%"
$ws.Range("D45").Value = "0.00006194"
$ws.Range("E45").Value = "-1.26%"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").Value = "0.06%"
$ws.Range("D47").Value = "0.8052"
$ws.Range("E47").Value = "-0.41%"
$ws.Range("D48").Value = "0.008388"
$ws.Range("E48").Value = "-16.05%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "0.06%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "0.13%"

# Restore the default (General) styling on the data rows of columns D and E
# so that only the cell text content changes and no residual "Text" number
# format remains.
$ws.Range($dataRange).Style = "Normal"
